$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.604.68'
$ws.Range('E2').Value = '  +7.00%  '
$ws.Range('D3').Value = '2.616.52'
$ws.Range('E3').Value = '  +7.36%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '184.60'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +14.21%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '581.03'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('E7').Value = '  -0.17%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.532'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +3.98%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.199'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +18.26%  '
$ws.Range('D10').Value = '2.616.63'
$ws.Range('E10').Value = '  +7.44%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.163'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.16%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.357'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +7.85%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.74'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +3.57%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.095.83'
$ws.Range('E14').Value = '  +7.33%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '73.502.37'
$ws.Range('E15').Value = '  +6.97%  '
$ws.Range('E16').Value = '  +6.17%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '26.08'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +12.27%  '
$ws.Range('D18').Value = '2.614.95'
$ws.Range('E18').Value = '  +7.38%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '9.09'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +31.04%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '11.81'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +12.06%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '370.49'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +9.25%  '
$ws.Range('E22').Value = '  +18.30%  '
$ws.Range('E23').Value = '  +6.48%  '
$ws.Range('E24').Value = '  +0.02%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '69.62'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +3.94%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '4.13'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +11.12%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.36'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +13.85%  '
$ws.Range('D28').Value = '2.737.10'
$ws.Range('E28').Value = '  +6.86%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '0.0₃0938'
$ws.Range('E30').Value = '  +14.50%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '523.35'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +22.32%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.39'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +19.80%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '7.59'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +6.40%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.73'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +8.49%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  +12.54%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '160.58'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.25%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '19.11'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +6.33%  '
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('E40').Value = '  -0.05%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '4.88'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +12.24%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.326'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +9.22%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.66'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +10.45%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '160.89'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +22.97%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.18'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +9.74%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.35'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +15.11%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '38.93'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +3.76%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0852'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +18.66%  '
$ws.Range('E49').Value = '  +8.36%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.527'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +9.73%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '20.61'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +22.34%  '
